# semana 19 de 2025
# Adds a new "19" column (column V) to the weekly IRA/UCI revision sheet,
# mirroring the same sparse fill pattern as the existing week columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for week 19, same style as the other week headers (bold/centered).
# Force text format first so the numeric-looking label "19" is stored as text,
# matching the other header cells (cod_pre..18).
$ws.Range("V1").NumberFormat = "@"
$ws.Range("V1").Value = "19"
$ws.Range("V1").Font.Bold = $true
$ws.Range("V1").HorizontalAlignment = -4108

# Data rows: most are 0, with a handful of non-zero counts.
$values = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    13 = 0
    14 = 0
    15 = 0
    17 = 0
    20 = 0
    23 = 0
    24 = 0
    26 = 0
    27 = 0
    28 = 5
    29 = 0
    30 = 0
    32 = 2
    33 = 0
    34 = 0
    35 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 22).Value = $values[$row]
}
